# Executive summary doc update
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1. Pain points header
Replace-Text `
  "The pain points/problems/needs/happiness (Vấn đề ,khó khăn khi thực hiện dự án) :" `
  "The pain points/problems/needs/happiness (Vấn đề) :"

# 2. First bullet
Replace-Text `
  "          Bất đồng ý kiến với đối tác" `
  "          Bất lợi trong việc phải phải tới rạp chiếu phim để mua vé , xem lịch chiếu phim"

# 3. Second bullet
Replace-Text `
  "          Thiếu nhân lực => không hoàn thành kịp hạn của dự án => đền bù hợp đồng" `
  "          Tốn thời gian xếp hàng để mua vé "

# 4. Third bullet
Replace-Text `
  "          Phân chia công việc chưa rõ ràng" `
  "          Sau khi mua vé xong thì phải chờ tới xuất chiếu"

# 5. Fourth bullet gets new text, and the bullet that followed it
#    ("Phát sinh thêm chi phí") is removed entirely.
Replace-Text `
  "          Trong quá trình thực hiện dự án nhà đầu tư giao yêu cầu mới và gấp rút" `
  "        Khác hàng tốn thời gian khi phải xem danh sách xuất chiếu (bình thường danh sách xuất chiếu hiển thị theo slide )"

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "         Phát sinh thêm chi phí") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# 6. Solutions bullet 1
Replace-Text `
  "          Kiểm soát chặt khâu nhân sự , báo cáo tình hình tiến trình dự án mỗi tuần " `
  "          Tạo web đặt vé xem phim ,người dùng có thể đặt vé trực tiếp mà không cần phải đi tới rạp  và xếp hàng"

# 7. Solutions bullet 2
Replace-Text `
  "          Có nhân lực dự trù để dễ dàng hỗ trợ nguồn nhân lực bị mất để đi đúng tiến độ" `
  "      Người dùng có thể xem các xuất chiếu mà không cần tới rạp chiếu phim"

# 8. Solutions bullet 3
Replace-Text `
  "          Team BA có kỹ năng tốt để đàm phán với nhà đầu tư để đẩy dự án đi vào quy trình bắt đầu " `
  " Chỉ cần đăt vé xem thông tin thời gian và tới đúng giờ là có thể xem phim liền ,người dùng không cần phải chờ "

# 9. Solutions bullet 4 gets new text, and the bullet that followed it
#    ("Kiểm soát chi phí ...") is removed entirely.
Replace-Text `
  "           Người quản lý nắm bắt rõ ràng yêu cầu mới của nhà đầu tư và nếu gấp rút thì thêm nhân lực dự trù để hỗ trợ cho team  " `
  "Người dùng thích thể loại phim nào thì có thể dễ dàng search thể loại phim "

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "            Kiểm soát chi phí nếu phát sinh phải biết phát sinh từ đâu và đưa giải pháp xử lý vi phạm đó") {
        $p.Range.Delete() | Out-Null
        break
    }
}
